$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Domain Name Health Diagnosis"
$ws.Range("B3").Value = "Comprehensive intellectualized domain name diagnosis to check the domain name health status at any time"
$ws.Range("B11").Value = "This domain name can be used normally"
$ws.Range("B12").Value = "The domain name will expire in ..."
$ws.Range("B17").Value = "The domain name resolution is normal"
$ws.Range("B19").Value = "Please enter correct domain name"
